# Sprint_Backlog.xlsx update
# - Fill in a "status" marker (column C) on the Backlog sheet for a batch
#   of existing rows: "u" (updated), "+" (added), "-" (removed).
# - Append a new backlog row ("Booking", priority 2, status "+").
# - Normalize the workbook's default cell style name.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Backlog")

$updatedRows = 19, 32, 34, 35, 42, 43
$addedRows   = 8, 14, 22, 27, 28, 30, 31, 33, 39, 41, 44
$removedRows = 36, 37, 38

foreach ($row in $updatedRows) {
    $ws.Cells.Item($row, 3).Value = "u"
}

foreach ($row in $addedRows) {
    $ws.Cells.Item($row, 3).Value = "+"
}

foreach ($row in $removedRows) {
    $ws.Cells.Item($row, 3).Value = "-"
}

# New backlog item: "Booking"
$ws.Cells.Item(45, 1).Value = 2
$ws.Cells.Item(45, 2).Value = "Booking"
$ws.Cells.Item(45, 3).Value = "+"

# Refresh the view: scroll position / zoom / selection on the Backlog sheet.
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 17
$win.ScrollColumn = 1
$win.Zoom = 100
$ws.Range("C44").Select()

# Default cell style rename (locale normalization: "Parasts" -> "Normal").
# Style.Name has no setter in the Excel object model, so recreate it.
$wb.Styles.Item("Parasts").Delete()
$wb.Styles.Add("Normal")
